# TestData.xlsx - "resolved issues of RW03"
#
# On the ScenarioMapping sheet, the SmokeTest ("D") flag is moved off the
# probateFormsRW02.feature scenarios and onto the probateFormsRW04.feature
# scenarios.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ScenarioMapping")

# probateFormsRW02.feature rows (TC_044 .. TC_074): SmokeTest Yes -> No
$ws.Range("D45:D75").Value = "No"

# probateFormsRW04.feature rows (TC_086 .. TC_097): SmokeTest No -> Yes
$ws.Range("D87:D98").Value = "Yes"

# Leave the cursor where the author left it when the sheet was saved.
$ws.Range("C98").Select() | Out-Null
